# feat: add 2022-Q1 data
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert a new worksheet "2022-Q1" right before the "总计" sheet,
#    in the same position/style as the other quarterly sheets.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$q1_2022 = $wb.Worksheets.Add($total)
$q1_2022.Name = "2022-Q1"

# Re-resolve the "总计" sheet by name: inserting a new sheet at its old
# slot rebinds the original $total handle to the newly inserted sheet.
$total = $wb.Worksheets.Item("总计")

# Headers (row 1)
$q1_2022.Range("B1").Value = "基金代码"
$q1_2022.Range("C1").Value = "基金名称"
$q1_2022.Range("D1").Value = "基金规模"
$q1_2022.Range("E1").Value = "股票总仓位"
$q1_2022.Range("F1").Value = "仓位占比"
$q1_2022.Range("G1").Value = "持有市值(亿元)"
$q1_2022.Range("H1").Value = "仓位排名"

# Data (row 2)
$q1_2022.Range("A2").Value = 0
$q1_2022.Range("B2").Value = "'513080"
$q1_2022.Range("C2").Value = "华安法国CAC40ETF（QDII）"
$q1_2022.Range("D2").Value = "'0.60"
$q1_2022.Range("E2").Value = "'96.69"
$q1_2022.Range("F2").Value = "'3.39"
$q1_2022.Range("G2").Value = "'0.0203"
$q1_2022.Range("H2").Value = 9

# Drop the quote-prefix styling picked up by forcing text above, so the
# cells end up with plain (unstyled) text, matching the other quarters.
$q1_2022.Range("B2:G2").ClearFormats()

# ---------------------------------------------------------------------
# 2. Update the "总计" (totals) sheet: prepend a new row for 2022-Q1,
#    shifting the existing date rows down and renumbering column A.
# ---------------------------------------------------------------------
$total.Rows.Item(2).Insert()
$total.Rows.Item(2).ClearFormats()

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.02

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q4"

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2021-Q3"

$total.Range("A5").Value = 3
$total.Range("B5").Value = "2021-Q1"

# A3 still carries the "index column" look (border + bold) from the
# original rows; reuse that same formatting for the new A2 cell.
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
